$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4 text content: plane launching/status terminology updated
$ws.Range("B4").Value = "land; status_grounded?; takeoff; status_flying?; "

# Update selection to B4
$ws.Range("B4").Select()
